# "Monte Alto added to MM347"
# Insert a new vessel row for "Monte Alto" into the MM347 (M347) project
# block, right before the existing "Maersk Sheerness" row (row 17),
# pushing all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 17 (shifts rows 17..225 down to 18..226).
$ws.Rows.Item(17).Insert()

# Populate the new row: NAME, IMO, CLASS, PROJECT
$ws.Cells.Item(17, 1).Value = "Monte Alto"
$ws.Cells.Item(17, 2).Value = 9283227
$ws.Cells.Item(17, 3).Value = "Monte"
$ws.Cells.Item(17, 4).Value = "M347"

# Reflect the author's last selection in the saved view state.
$ws.Range("F16").Select() | Out-Null
